# Append 10 new log rows (137-146) to the "check_availability" sheet.
# Note: column E values look like plain dates ("2024-09-29"); a leading
# apostrophe is used so Excel stores them as literal text instead of
# auto-converting them to date serial numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(137, 1).Value = '2024-09-29 03:54:47'
$ws.Cells.Item(137, 2).Value = 'check_availability'
$ws.Cells.Item(137, 3).Value = 'https://example.com'
$ws.Cells.Item(137, 4).Value = 'Checked availability: Selected or default date current date is available for booking.'
$ws.Cells.Item(137, 5).Value = "'2024-09-29"
$ws.Cells.Item(137, 6).Value = '03:54:47'

$ws.Cells.Item(138, 1).Value = '2024-09-29 03:54:47'
$ws.Cells.Item(138, 2).Value = 'check_availability'
$ws.Cells.Item(138, 3).Value = 'https://example.com'
$ws.Cells.Item(138, 4).Value = 'Failed to check availability: Failed to check availability'
$ws.Cells.Item(138, 5).Value = "'2024-09-29"
$ws.Cells.Item(138, 6).Value = '03:54:47'

$ws.Cells.Item(139, 1).Value = '2024-09-29 03:54:47'
$ws.Cells.Item(139, 2).Value = 'check_availability'
$ws.Cells.Item(139, 3).Value = 'https://example.com'
$ws.Cells.Item(139, 4).Value = 'Checked availability: No availability for the selected date.'
$ws.Cells.Item(139, 5).Value = "'2024-09-29"
$ws.Cells.Item(139, 6).Value = '03:54:47'

$ws.Cells.Item(140, 1).Value = '2024-09-29 03:54:48'
$ws.Cells.Item(140, 2).Value = 'check_availability'
$ws.Cells.Item(140, 3).Value = 'https://example.com'
$ws.Cells.Item(140, 4).Value = 'Checked availability: Selected or default date is available for booking.'
$ws.Cells.Item(140, 5).Value = "'2024-09-29"
$ws.Cells.Item(140, 6).Value = '03:54:48'

$ws.Cells.Item(141, 1).Value = '2024-09-29 03:54:49'
$ws.Cells.Item(141, 2).Value = 'check_availability'
$ws.Cells.Item(141, 3).Value = 'https://example.com'
$ws.Cells.Item(141, 4).Value = 'Failed to check availability: Failed to check availability'
$ws.Cells.Item(141, 5).Value = "'2024-09-29"
$ws.Cells.Item(141, 6).Value = '03:54:49'

$ws.Cells.Item(142, 1).Value = '2024-09-29 03:56:04'
$ws.Cells.Item(142, 2).Value = 'check_availability'
$ws.Cells.Item(142, 3).Value = 'https://example.com'
$ws.Cells.Item(142, 4).Value = 'Checked availability: Selected or default date current date is available for booking.'
$ws.Cells.Item(142, 5).Value = "'2024-09-29"
$ws.Cells.Item(142, 6).Value = '03:56:04'

$ws.Cells.Item(143, 1).Value = '2024-09-29 03:56:05'
$ws.Cells.Item(143, 2).Value = 'check_availability'
$ws.Cells.Item(143, 3).Value = 'https://example.com'
$ws.Cells.Item(143, 4).Value = 'Failed to check availability: Failed to check availability'
$ws.Cells.Item(143, 5).Value = "'2024-09-29"
$ws.Cells.Item(143, 6).Value = '03:56:05'

$ws.Cells.Item(144, 1).Value = '2024-09-29 03:56:05'
$ws.Cells.Item(144, 2).Value = 'check_availability'
$ws.Cells.Item(144, 3).Value = 'https://example.com'
$ws.Cells.Item(144, 4).Value = 'Checked availability: No availability for the selected date.'
$ws.Cells.Item(144, 5).Value = "'2024-09-29"
$ws.Cells.Item(144, 6).Value = '03:56:05'

$ws.Cells.Item(145, 1).Value = '2024-09-29 03:56:05'
$ws.Cells.Item(145, 2).Value = 'check_availability'
$ws.Cells.Item(145, 3).Value = 'https://example.com'
$ws.Cells.Item(145, 4).Value = 'Checked availability: Selected or default date is available for booking.'
$ws.Cells.Item(145, 5).Value = "'2024-09-29"
$ws.Cells.Item(145, 6).Value = '03:56:05'

$ws.Cells.Item(146, 1).Value = '2024-09-29 03:56:06'
$ws.Cells.Item(146, 2).Value = 'check_availability'
$ws.Cells.Item(146, 3).Value = 'https://example.com'
$ws.Cells.Item(146, 4).Value = 'Failed to check availability: Failed to check availability'
$ws.Cells.Item(146, 5).Value = "'2024-09-29"
$ws.Cells.Item(146, 6).Value = '03:56:06'
